$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format to avoid numeric auto-conversion while we set string values
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.307.66"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "1.588.97"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "19.44"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "1.811.76"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "1.578.28"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "64.37"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "26.316.00"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "7.50"
$ws.Range("E19").Value = "  +6.13%  "
$ws.Range("D20").Value = "211.19"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "4.27"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "8.95"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").Value = "144.45"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "7.04"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "0.0506"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "1.310.58"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").Value = "0.612"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  -9.67%  "
$ws.Range("D40").Value = "0.806"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "5.61"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").Value = "0.768"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "2.13"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "62.30"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "1.724.33"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").Value = "87.47"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("E48").Value = "  -5.21%  "
$ws.Range("D49").Value = "0.0505"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("E50").Value = "  -4.88%  "
$ws.Range("E51").Value = "  -0.45%  "

# Restore default formatting on column D so styles match the original (plain inline/shared text, no style id)
$priceRange.ClearFormats()

